$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.022.73"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.547.80"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "197.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.654"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000304"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "4.117.43"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "605.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "70.194.24"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "3.555.53"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +22.44%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "0.0₃0839"
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("D36").Value = "3.766.12"
$ws.Range("E36").Value = "  +7.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "487.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.35%  "
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "

# Row 51: Mantle -> Monero (full row change)
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
